$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule row (was "Good Morning").
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the last user selection being on E8.
$ws.Range("E8").Select()
